$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 - The Little Girl
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "TheLittleGirl"
$ws.Range("C50").Value = "The Little Girl"
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = "Labyrinth"
$ws.Range("I50").Value = 224
$ws.Range("J50").Value = 2464

# Row 51 - The Lady in White
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "TheLadyinWhite"
$ws.Range("C51").Value = "The Lady in White"
$ws.Range("D51").Value = 2
$ws.Range("E51").Value = "Labyrinth"
$ws.Range("I51").Value = 48
$ws.Range("J51").Value = 1040
